# Appends 27 new transaction rows (rows 38-64) for the week ending 2021-01-17
# to the "Konto" sheet, following the existing Datum/Receipt Number/Konto/
# Beskrivning/Debet/Kredit layout. New dimension becomes A1:F64.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-use the same date/time number format already applied to column A (style index 2,
# numFmtId 165 "YYYY-MM-DD HH:MM:SS") so new cells match existing Datum cells exactly.
$dateFormat = $ws.Cells.Item(37, 1).NumberFormat

# Row 38
$ws.Cells.Item(38, 1).Value = 44207
$ws.Cells.Item(38, 1).NumberFormat = $dateFormat
$ws.Cells.Item(38, 3).Value = 4010
$ws.Cells.Item(38, 4).Value = "SvE*LINDH BROS AB K0135"
$ws.Cells.Item(38, 5).Value = 788.35

# Row 39
$ws.Cells.Item(39, 1).Value = 44207
$ws.Cells.Item(39, 1).NumberFormat = $dateFormat
$ws.Cells.Item(39, 3).Value = 2641
$ws.Cells.Item(39, 4).Value = "SvE*LINDH BROS AB K0135"
$ws.Cells.Item(39, 5).Value = 197.09

# Row 40
$ws.Cells.Item(40, 1).Value = 44207
$ws.Cells.Item(40, 1).NumberFormat = $dateFormat
$ws.Cells.Item(40, 3).Value = 1930
$ws.Cells.Item(40, 4).Value = "SvE*LINDH BROS AB K0135"
$ws.Cells.Item(40, 6).Value = 985.44

# Row 41
$ws.Cells.Item(41, 1).Value = 44208
$ws.Cells.Item(41, 1).NumberFormat = $dateFormat
$ws.Cells.Item(41, 2).Value = 4121803
$ws.Cells.Item(41, 3).Value = 3011
$ws.Cells.Item(41, 4).Value = "Order 4121803 Swish +46723656673"
$ws.Cells.Item(41, 6).Value = 806.25

# Row 42
$ws.Cells.Item(42, 1).Value = 44208
$ws.Cells.Item(42, 1).NumberFormat = $dateFormat
$ws.Cells.Item(42, 2).Value = 4121803
$ws.Cells.Item(42, 3).Value = 2611
$ws.Cells.Item(42, 4).Value = "Order 4121803 Swish +46723656673"
$ws.Cells.Item(42, 6).Value = 96.75

# Row 43
$ws.Cells.Item(43, 1).Value = 44208
$ws.Cells.Item(43, 1).NumberFormat = $dateFormat
$ws.Cells.Item(43, 2).Value = 4121803
$ws.Cells.Item(43, 3).Value = 1930
$ws.Cells.Item(43, 4).Value = "Order 4121803 Swish +46723656673"
$ws.Cells.Item(43, 5).Value = 903

# Row 44
$ws.Cells.Item(44, 1).Value = 44209
$ws.Cells.Item(44, 1).NumberFormat = $dateFormat
$ws.Cells.Item(44, 2).Value = 6130504
$ws.Cells.Item(44, 3).Value = 3011
$ws.Cells.Item(44, 4).Value = "Order 6130504 Swish +46709526084"
$ws.Cells.Item(44, 6).Value = 1242.86

# Row 45
$ws.Cells.Item(45, 1).Value = 44209
$ws.Cells.Item(45, 1).NumberFormat = $dateFormat
$ws.Cells.Item(45, 2).Value = 6130504
$ws.Cells.Item(45, 3).Value = 2611
$ws.Cells.Item(45, 4).Value = "Order 6130504 Swish +46709526084"
$ws.Cells.Item(45, 6).Value = 149.14

# Row 46
$ws.Cells.Item(46, 1).Value = 44209
$ws.Cells.Item(46, 1).NumberFormat = $dateFormat
$ws.Cells.Item(46, 2).Value = 6130504
$ws.Cells.Item(46, 3).Value = 1930
$ws.Cells.Item(46, 4).Value = "Order 6130504 Swish +46709526084"
$ws.Cells.Item(46, 5).Value = 1392

# Row 47
$ws.Cells.Item(47, 1).Value = 44209
$ws.Cells.Item(47, 1).NumberFormat = $dateFormat
$ws.Cells.Item(47, 2).Value = 3131821
$ws.Cells.Item(47, 3).Value = 3011
$ws.Cells.Item(47, 4).Value = "Order 3131821 Swish +46738066249"
$ws.Cells.Item(47, 6).Value = 398.21

# Row 48
$ws.Cells.Item(48, 1).Value = 44209
$ws.Cells.Item(48, 1).NumberFormat = $dateFormat
$ws.Cells.Item(48, 2).Value = 3131821
$ws.Cells.Item(48, 3).Value = 2611
$ws.Cells.Item(48, 4).Value = "Order 3131821 Swish +46738066249"
$ws.Cells.Item(48, 6).Value = 47.79

# Row 49
$ws.Cells.Item(49, 1).Value = 44209
$ws.Cells.Item(49, 1).NumberFormat = $dateFormat
$ws.Cells.Item(49, 2).Value = 3131821
$ws.Cells.Item(49, 3).Value = 1930
$ws.Cells.Item(49, 4).Value = "Order 3131821 Swish +46738066249"
$ws.Cells.Item(49, 5).Value = 446

# Row 50
$ws.Cells.Item(50, 1).Value = 44209
$ws.Cells.Item(50, 1).NumberFormat = $dateFormat
$ws.Cells.Item(50, 2).Value = 2131916
$ws.Cells.Item(50, 3).Value = 3011
$ws.Cells.Item(50, 4).Value = "Order 2131916 Swish +46702597315"
$ws.Cells.Item(50, 6).Value = 502.68

# Row 51
$ws.Cells.Item(51, 1).Value = 44209
$ws.Cells.Item(51, 1).NumberFormat = $dateFormat
$ws.Cells.Item(51, 2).Value = 2131916
$ws.Cells.Item(51, 3).Value = 2611
$ws.Cells.Item(51, 4).Value = "Order 2131916 Swish +46702597315"
$ws.Cells.Item(51, 6).Value = 60.32

# Row 52
$ws.Cells.Item(52, 1).Value = 44209
$ws.Cells.Item(52, 1).NumberFormat = $dateFormat
$ws.Cells.Item(52, 2).Value = 2131916
$ws.Cells.Item(52, 3).Value = 1930
$ws.Cells.Item(52, 4).Value = "Order 2131916 Swish +46702597315"
$ws.Cells.Item(52, 5).Value = 563

# Row 53
$ws.Cells.Item(53, 1).Value = 44210
$ws.Cells.Item(53, 1).NumberFormat = $dateFormat
$ws.Cells.Item(53, 3).Value = 4010
$ws.Cells.Item(53, 4).Value = "MATVA.RLDEN VA.LLINGBY K6885"
$ws.Cells.Item(53, 5).Value = 30.84

# Row 54
$ws.Cells.Item(54, 1).Value = 44210
$ws.Cells.Item(54, 1).NumberFormat = $dateFormat
$ws.Cells.Item(54, 3).Value = 2645
$ws.Cells.Item(54, 4).Value = "MATVA.RLDEN VA.LLINGBY K6885"
$ws.Cells.Item(54, 5).Value = 3.7

# Row 55
$ws.Cells.Item(55, 1).Value = 44210
$ws.Cells.Item(55, 1).NumberFormat = $dateFormat
$ws.Cells.Item(55, 3).Value = 1930
$ws.Cells.Item(55, 4).Value = "MATVA.RLDEN VA.LLINGBY K6885"
$ws.Cells.Item(55, 6).Value = 34.54

# Row 56
$ws.Cells.Item(56, 1).Value = 44211
$ws.Cells.Item(56, 1).NumberFormat = $dateFormat
$ws.Cells.Item(56, 2).Value = 3152209
$ws.Cells.Item(56, 3).Value = 3011
$ws.Cells.Item(56, 4).Value = "Order 3152209 Swish +46722017122"
$ws.Cells.Item(56, 6).Value = 398.21

# Row 57
$ws.Cells.Item(57, 1).Value = 44211
$ws.Cells.Item(57, 1).NumberFormat = $dateFormat
$ws.Cells.Item(57, 2).Value = 3152209
$ws.Cells.Item(57, 3).Value = 2611
$ws.Cells.Item(57, 4).Value = "Order 3152209 Swish +46722017122"
$ws.Cells.Item(57, 6).Value = 47.79

# Row 58
$ws.Cells.Item(58, 1).Value = 44211
$ws.Cells.Item(58, 1).NumberFormat = $dateFormat
$ws.Cells.Item(58, 2).Value = 3152209
$ws.Cells.Item(58, 3).Value = 1930
$ws.Cells.Item(58, 4).Value = "Order 3152209 Swish +46722017122"
$ws.Cells.Item(58, 5).Value = 446

# Row 59
$ws.Cells.Item(59, 1).Value = 44212
$ws.Cells.Item(59, 1).NumberFormat = $dateFormat
$ws.Cells.Item(59, 3).Value = 4010
$ws.Cells.Item(59, 4).Value = "FRESH STOCKHOLM HÄSSEL K6885"
$ws.Cells.Item(59, 5).Value = 172.66

# Row 60
$ws.Cells.Item(60, 1).Value = 44212
$ws.Cells.Item(60, 1).NumberFormat = $dateFormat
$ws.Cells.Item(60, 3).Value = 2645
$ws.Cells.Item(60, 4).Value = "FRESH STOCKHOLM HÄSSEL K6885"
$ws.Cells.Item(60, 5).Value = 20.72

# Row 61
$ws.Cells.Item(61, 1).Value = 44212
$ws.Cells.Item(61, 1).NumberFormat = $dateFormat
$ws.Cells.Item(61, 3).Value = 1930
$ws.Cells.Item(61, 4).Value = "FRESH STOCKHOLM HÄSSEL K6885"
$ws.Cells.Item(61, 6).Value = 193.38

# Row 62
$ws.Cells.Item(62, 1).Value = 44213
$ws.Cells.Item(62, 1).NumberFormat = $dateFormat
$ws.Cells.Item(62, 3).Value = 4010
$ws.Cells.Item(62, 4).Value = "MATVARLDEN VEDD K0135"
$ws.Cells.Item(62, 5).Value = 442.54

# Row 63
$ws.Cells.Item(63, 1).Value = 44213
$ws.Cells.Item(63, 1).NumberFormat = $dateFormat
$ws.Cells.Item(63, 3).Value = 2645
$ws.Cells.Item(63, 4).Value = "MATVARLDEN VEDD K0135"
$ws.Cells.Item(63, 5).Value = 53.1

# Row 64
$ws.Cells.Item(64, 1).Value = 44213
$ws.Cells.Item(64, 1).NumberFormat = $dateFormat
$ws.Cells.Item(64, 3).Value = 1930
$ws.Cells.Item(64, 4).Value = "MATVARLDEN VEDD K0135"
$ws.Cells.Item(64, 6).Value = 495.64
